# Update recomputed TPM-derived NATMI edge statistics for Hspg2-Itga2 (YoungD2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = "204.558136"
$ws.Range("H2").Value = "613.674408"
$ws.Range("I2").Value = "0.60178627893129"
$ws.Range("J2").Value = "0.6017862789312901"
$ws.Range("M2").Value = "3.339352"
$ws.Range("N2").Value = "10.018056"
$ws.Range("O2").Value = "0.6054960700393903"
$ws.Range("P2").Value = "0.6054960700393903"
$ws.Range("Q2").Value = "683.091620567872"
$ws.Range("R2").Value = "6147.824585110849"
$ws.Range("S2").Value = "0.3643792268965245"
$ws.Range("T2").Value = "0.3643792268965245"
# Row 3
$ws.Range("G3").Value = "204.558136"
$ws.Range("H3").Value = "613.674408"
$ws.Range("I3").Value = "0.60178627893129"
$ws.Range("J3").Value = "0.6017862789312901"
$ws.Range("O3").Value = "0.2540955070726236"
$ws.Range("P3").Value = "0.2540955070726236"
$ws.Range("Q3").Value = "286.658362116144"
$ws.Range("R3").Value = "2579.925259045296"
$ws.Range("S3").Value = "0.1529111896943935"
$ws.Range("T3").Value = "0.1529111896943934"
# Row 4
$ws.Range("G4").Value = "204.558136"
$ws.Range("H4").Value = "613.674408"
$ws.Range("I4").Value = "0.60178627893129"
$ws.Range("J4").Value = "0.6017862789312901"
$ws.Range("K4").Value = "2"
$ws.Range("L4").Value = "0.6666666666666666"
$ws.Range("M4").Value = "0.1338136666666667"
$ws.Range("N4").Value = "0.401441"
$ws.Range("O4").Value = "0.02426328499787613"
$ws.Range("P4").Value = "0.02426328499787612"
$ws.Range("Q4").Value = "27.37267422465867"
$ws.Range("R4").Value = "246.354068021928"
$ws.Range("S4").Value = "0.01460131199352127"
$ws.Range("T4").Value = "0.01460131199352127"
# Row 5
$ws.Range("G5").Value = "204.558136"
$ws.Range("H5").Value = "613.674408"
$ws.Range("I5").Value = "0.60178627893129"
$ws.Range("J5").Value = "0.6017862789312901"
$ws.Range("M5").Value = "0.6405483333333334"
$ws.Range("N5").Value = "1.921645"
$ws.Range("O5").Value = "0.11614513789011"
$ws.Range("P5").Value = "0.11614513789011"
$ws.Range("Q5").Value = "131.0293730845733"
$ws.Range("R5").Value = "1179.26435776116"
$ws.Range("S5").Value = "0.06989455034685089"
$ws.Range("T5").Value = "0.06989455034685089"
# Row 6
$ws.Range("I6").Value = "0.3090998990957371"
$ws.Range("J6").Value = "0.3090998990957372"
$ws.Range("M6").Value = "3.339352"
$ws.Range("N6").Value = "10.018056"
$ws.Range("O6").Value = "0.6054960700393903"
$ws.Range("P6").Value = "0.6054960700393903"
$ws.Range("Q6").Value = "350.8613578987574"
$ws.Range("R6").Value = "3157.752221088817"
$ws.Range("S6").Value = "0.1871587741520409"
$ws.Range("T6").Value = "0.1871587741520409"
# Row 7
$ws.Range("I7").Value = "0.3090998990957371"
$ws.Range("J7").Value = "0.3090998990957372"
$ws.Range("O7").Value = "0.2540955070726236"
$ws.Range("P7").Value = "0.2540955070726236"
$ws.Range("S7").Value = "0.07854089559682811"
$ws.Range("T7").Value = "0.07854089559682811"
# Row 8
$ws.Range("I8").Value = "0.3090998990957371"
$ws.Range("J8").Value = "0.3090998990957372"
$ws.Range("K8").Value = "2"
$ws.Range("L8").Value = "0.6666666666666666"
$ws.Range("M8").Value = "0.1338136666666667"
$ws.Range("N8").Value = "0.401441"
$ws.Range("O8").Value = "0.02426328499787613"
$ws.Range("P8").Value = "0.02426328499787612"
$ws.Range("Q8").Value = "14.05962737443623"
$ws.Range("R8").Value = "126.536646369926"
$ws.Range("S8").Value = "0.007499778944574622"
$ws.Range("T8").Value = "0.007499778944574622"
# Row 9
$ws.Range("I9").Value = "0.3090998990957371"
$ws.Range("J9").Value = "0.3090998990957372"
$ws.Range("M9").Value = "0.6405483333333334"
$ws.Range("N9").Value = "1.921645"
$ws.Range("O9").Value = "0.11614513789011"
$ws.Range("P9").Value = "0.11614513789011"
$ws.Range("Q9").Value = "67.30157768127445"
$ws.Range("R9").Value = "605.7141991314701"
$ws.Range("S9").Value = "0.03590045040229348"
$ws.Range("T9").Value = "0.03590045040229348"
# Row 10
$ws.Range("G10").Value = "0.1651866666666667"
$ws.Range("H10").Value = "0.49556"
$ws.Range("I10").Value = "0.0004859599887163456"
$ws.Range("J10").Value = "0.0004859599887163457"
$ws.Range("M10").Value = "3.339352"
$ws.Range("N10").Value = "10.018056"
$ws.Range("O10").Value = "0.6054960700393903"
$ws.Range("P10").Value = "0.6054960700393903"
$ws.Range("Q10").Value = "0.5516164257066668"
$ws.Range("R10").Value = "4.964547831360001"
$ws.Range("S10").Value = "0.0002942468633641337"
$ws.Range("T10").Value = "0.0002942468633641338"
# Row 11
$ws.Range("G11").Value = "0.1651866666666667"
$ws.Range("H11").Value = "0.49556"
$ws.Range("I11").Value = "0.0004859599887163456"
$ws.Range("J11").Value = "0.0004859599887163457"
$ws.Range("O11").Value = "0.2540955070726236"
$ws.Range("P11").Value = "0.2540955070726236"
$ws.Range("Q11").Value = "0.2314849960800001"
$ws.Range("R11").Value = "2.08336496472"
$ws.Range("S11").Value = "0.0001234802497498863"
$ws.Range("T11").Value = "0.0001234802497498863"
# Row 12
$ws.Range("G12").Value = "0.1651866666666667"
$ws.Range("H12").Value = "0.49556"
$ws.Range("I12").Value = "0.0004859599887163456"
$ws.Range("J12").Value = "0.0004859599887163457"
$ws.Range("K12").Value = "2"
$ws.Range("L12").Value = "0.6666666666666666"
$ws.Range("M12").Value = "0.1338136666666667"
$ws.Range("N12").Value = "0.401441"
$ws.Range("O12").Value = "0.02426328499787613"
$ws.Range("P12").Value = "0.02426328499787612"
$ws.Range("Q12").Value = "0.02210423355111112"
$ws.Range("R12").Value = "0.19893810196"
$ws.Range("S12").Value = "0.00001179098570378936"
$ws.Range("T12").Value = "0.00001179098570378936"
# Row 13
$ws.Range("G13").Value = "0.1651866666666667"
$ws.Range("H13").Value = "0.49556"
$ws.Range("I13").Value = "0.0004859599887163456"
$ws.Range("J13").Value = "0.0004859599887163457"
$ws.Range("M13").Value = "0.6405483333333334"
$ws.Range("N13").Value = "1.921645"
$ws.Range("O13").Value = "0.11614513789011"
$ws.Range("P13").Value = "0.11614513789011"
$ws.Range("Q13").Value = "0.1058100440222222"
$ws.Range("R13").Value = "0.9522903962"
$ws.Range("S13").Value = "0.00005644188989853628"
$ws.Range("T13").Value = "0.00005644188989853628"
# Row 14
$ws.Range("G14").Value = "29.46746633333333"
$ws.Range("H14").Value = "88.402399"
$ws.Range("I14").Value = "0.08668986363011115"
$ws.Range("J14").Value = "0.08668986363011116"
$ws.Range("M14").Value = "3.339352"
$ws.Range("N14").Value = "10.018056"
$ws.Range("O14").Value = "0.6054960700393903"
$ws.Range("P14").Value = "0.6054960700393903"
$ws.Range("Q14").Value = "98.40224263514935"
$ws.Range("R14").Value = "885.6201837163442"
$ws.Range("S14").Value = "0.05249037174028297"
$ws.Range("T14").Value = "0.05249037174028298"
# Row 15
$ws.Range("G15").Value = "29.46746633333333"
$ws.Range("H15").Value = "88.402399"
$ws.Range("I15").Value = "0.08668986363011115"
$ws.Range("J15").Value = "0.08668986363011116"
$ws.Range("O15").Value = "0.2540955070726236"
$ws.Range("P15").Value = "0.2540955070726236"
$ws.Range("Q15").Value = "41.29435181608201"
$ws.Range("R15").Value = "371.649166344738"
$ws.Range("S15").Value = "0.02202750485714968"
$ws.Range("T15").Value = "0.02202750485714968"
# Row 16
$ws.Range("G16").Value = "29.46746633333333"
$ws.Range("H16").Value = "88.402399"
$ws.Range("I16").Value = "0.08668986363011115"
$ws.Range("J16").Value = "0.08668986363011116"
$ws.Range("K16").Value = "2"
$ws.Range("L16").Value = "0.6666666666666666"
$ws.Range("M16").Value = "0.1338136666666667"
$ws.Range("N16").Value = "0.401441"
$ws.Range("O16").Value = "0.02426328499787613"
$ws.Range("P16").Value = "0.02426328499787612"
$ws.Range("Q16").Value = "3.94314971743989"
$ws.Range("R16").Value = "35.488347456959"
$ws.Range("S16").Value = "0.002103380867684403"
$ws.Range("T16").Value = "0.002103380867684403"
# Row 17
$ws.Range("G17").Value = "29.46746633333333"
$ws.Range("H17").Value = "88.402399"
$ws.Range("I17").Value = "0.08668986363011115"
$ws.Range("J17").Value = "0.08668986363011116"
$ws.Range("M17").Value = "0.6405483333333334"
$ws.Range("N17").Value = "1.921645"
$ws.Range("O17").Value = "0.11614513789011"
$ws.Range("P17").Value = "0.11614513789011"
$ws.Range("Q17").Value = "18.87533644737278"
$ws.Range("R17").Value = "169.878028026355"
$ws.Range("S17").Value = "0.01006860616499409"
$ws.Range("T17").Value = "0.0100686061649941"
# Row 18
$ws.Range("G18").Value = "0.54608"
$ws.Range("H18").Value = "1.63824"
$ws.Range("I18").Value = "0.001606503938805929"
$ws.Range("J18").Value = "0.001606503938805929"
$ws.Range("M18").Value = "3.339352"
$ws.Range("N18").Value = "10.018056"
$ws.Range("O18").Value = "0.6054960700393903"
$ws.Range("P18").Value = "0.6054960700393903"
$ws.Range("Q18").Value = "1.82355334016"
$ws.Range("R18").Value = "16.41198006144"
$ws.Range("S18").Value = "0.000972731821449791"
$ws.Range("T18").Value = "0.0009727318214497911"
# Row 19
$ws.Range("G19").Value = "0.54608"
$ws.Range("H19").Value = "1.63824"
$ws.Range("I19").Value = "0.001606503938805929"
$ws.Range("J19").Value = "0.001606503938805929"
$ws.Range("O19").Value = "0.2540955070726236"
$ws.Range("P19").Value = "0.2540955070726236"
$ws.Range("Q19").Value = "0.7652513923200002"
$ws.Range("R19").Value = "6.88726253088"
$ws.Range("S19").Value = "0.0004082054329450595"
$ws.Range("T19").Value = "0.0004082054329450595"
# Row 20
$ws.Range("G20").Value = "0.54608"
$ws.Range("H20").Value = "1.63824"
$ws.Range("I20").Value = "0.001606503938805929"
$ws.Range("J20").Value = "0.001606503938805929"
$ws.Range("K20").Value = "2"
$ws.Range("L20").Value = "0.6666666666666666"
$ws.Range("M20").Value = "0.1338136666666667"
$ws.Range("N20").Value = "0.401441"
$ws.Range("O20").Value = "0.02426328499787613"
$ws.Range("P20").Value = "0.02426328499787612"
$ws.Range("Q20").Value = "0.07307296709333334"
$ws.Range("R20").Value = "0.6576567038400001"
$ws.Range("S20").Value = "0.00003897906291745879"
$ws.Range("T20").Value = "0.00003897906291745879"
# Row 21
$ws.Range("G21").Value = "0.54608"
$ws.Range("H21").Value = "1.63824"
$ws.Range("I21").Value = "0.001606503938805929"
$ws.Range("J21").Value = "0.001606503938805929"
$ws.Range("M21").Value = "0.6405483333333334"
$ws.Range("N21").Value = "1.921645"
$ws.Range("O21").Value = "0.11614513789011"
$ws.Range("P21").Value = "0.11614513789011"
$ws.Range("Q21").Value = "0.3497906338666667"
$ws.Range("R21").Value = "3.1481157048"
$ws.Range("S21").Value = "0.0001865876214936195"
$ws.Range("T21").Value = "0.0001865876214936195"
# Row 22
$ws.Range("G22").Value = "0.112681"
$ws.Range("H22").Value = "0.338043"
$ws.Range("I22").Value = "0.000331494415339494"
$ws.Range("J22").Value = "0.000331494415339494"
$ws.Range("M22").Value = "3.339352"
$ws.Range("N22").Value = "10.018056"
$ws.Range("O22").Value = "0.6054960700393903"
$ws.Range("P22").Value = "0.6054960700393903"
$ws.Range("Q22").Value = "0.3762815227120001"
$ws.Range("R22").Value = "3.386533704408001"
$ws.Range("S22").Value = "0.000200718565728069"
$ws.Range("T22").Value = "0.000200718565728069"
# Row 23
$ws.Range("G23").Value = "0.112681"
$ws.Range("H23").Value = "0.338043"
$ws.Range("I23").Value = "0.000331494415339494"
$ws.Range("J23").Value = "0.000331494415339494"
$ws.Range("O23").Value = "0.2540955070726236"
$ws.Range("P23").Value = "0.2540955070726236"
$ws.Range("Q23").Value = "0.157905970074"
$ws.Range("R23").Value = "1.421153730666"
$ws.Range("S23").Value = "0.00008423124155743162"
$ws.Range("T23").Value = "0.00008423124155743162"
# Row 24
$ws.Range("G24").Value = "0.112681"
$ws.Range("H24").Value = "0.338043"
$ws.Range("I24").Value = "0.000331494415339494"
$ws.Range("J24").Value = "0.000331494415339494"
$ws.Range("K24").Value = "2"
$ws.Range("L24").Value = "0.6666666666666666"
$ws.Range("M24").Value = "0.1338136666666667"
$ws.Range("N24").Value = "0.401441"
$ws.Range("O24").Value = "0.02426328499787613"
$ws.Range("P24").Value = "0.02426328499787612"
$ws.Range("Q24").Value = "0.01507825777366667"
$ws.Range("R24").Value = "0.135704319963"
$ws.Range("S24").Value = "0.000008043143474586462"
$ws.Range("T24").Value = "0.000008043143474586462"
# Row 25
$ws.Range("G25").Value = "0.112681"
$ws.Range("H25").Value = "0.338043"
$ws.Range("I25").Value = "0.000331494415339494"
$ws.Range("J25").Value = "0.000331494415339494"
$ws.Range("M25").Value = "0.6405483333333334"
$ws.Range("N25").Value = "1.921645"
$ws.Range("O25").Value = "0.11614513789011"
$ws.Range("P25").Value = "0.11614513789011"
$ws.Range("Q25").Value = "0.07217762674833335"
$ws.Range("R25").Value = "0.6495986407350001"
$ws.Range("S25").Value = "0.00003850146457940694"
$ws.Range("T25").Value = "0.00003850146457940694"
